$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "45.857.61"
Set-TextValue "E2" "  +0.87%  "
Set-TextValue "D3" "2.516.96"
Set-TextValue "E3" "  +8.02%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.13%  "
Set-TextValue "D5" "302.96"
Set-TextValue "E5" "  +1.44%  "
Set-TextValue "D6" "100.62"
Set-TextValue "E6" "  +4.10%  "
Set-TextValue "D7" "0.596"
Set-TextValue "E7" "  +5.04%  "
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D9" "0.557"
Set-TextValue "E9" "  +9.71%  "
Set-TextValue "D10" "38.30"
Set-TextValue "E10" "  +12.58%  "
Set-TextValue "D11" "0.0815"
Set-TextValue "E11" "  +4.79%  "
Set-TextValue "D12" "7.83"
Set-TextValue "E12" "  +11.48%  "
Set-TextValue "D13" "2.913.34"
Set-TextValue "E13" "  +8.63%  "
Set-TextValue "E14" "  +2.65%  "
Set-TextValue "D15" "2.526.29"
Set-TextValue "E15" "  +7.26%  "
Set-TextValue "D16" "0.884"
Set-TextValue "E16" "  +10.88%  "
Set-TextValue "D17" "14.72"
Set-TextValue "E17" "  +8.65%  "
Set-TextValue "D18" "45.959.90"
Set-TextValue "E18" "  +1.03%  "
Set-TextValue "D19" "13.59"
Set-TextValue "E19" "  +9.59%  "
Set-TextValue "D20" "0.0₃0984"
Set-TextValue "E20" "  +2.32%  "
Set-TextValue "D21" "6.56"
Set-TextValue "E21" "  +11.83%  "
Set-TextValue "D22" "69.39"
Set-TextValue "E22" "  +5.83%  "
Set-TextValue "D23" "251.06"
Set-TextValue "E23" "  +3.88%  "
Set-TextValue "D24" "2.93"
Set-TextValue "E24" "  +5.32%  "
Set-TextValue "D25" "2.08"
Set-TextValue "E25" "  +11.54%  "
Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  -0.20%  "
Set-TextValue "D27" "41.11"
Set-TextValue "E27" "  +3.24%  "
Set-TextValue "D28" "23.64"
Set-TextValue "E28" "  +15.97%  "
Set-TextValue "D29" "10.26"
Set-TextValue "E29" "  +7.35%  "
Set-TextValue "D30" "2.25"
Set-TextValue "E30" "  +1.62%  "
Set-TextValue "D31" "3.80"
Set-TextValue "E31" "  +7.53%  "
Set-TextValue "D32" "5.92"
Set-TextValue "E32" "  +11.25%  "
Set-TextValue "D33" "2.89"
Set-TextValue "E33" "  +6.86%  "
Set-TextValue "B34" "Hedera"
Set-TextValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.0823"
Set-TextValue "E34" "  +7.71%  "
Set-TextValue "B35" "ARBITRUM"
Set-TextValue "C35" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "2.16"
Set-TextValue "E35" "  +23.59%  "
Set-TextValue "B36" "Monero"
Set-TextValue "C36" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D36" "149.21"
Set-TextValue "E36" "  +3.29%  "
Set-TextValue "E37" "  +6.34%  "
Set-TextValue "D38" "0.120"
Set-TextValue "E38" "  +4.05%  "
Set-TextValue "D39" "16.12"
Set-TextValue "E39" "  +5.86%  "
Set-TextValue "D40" "4.15"
Set-TextValue "E40" "  +8.63%  "
Set-TextValue "D41" "0.0321"
Set-TextValue "E41" "  +9.31%  "
Set-TextValue "D42" "3.46"
Set-TextValue "E42" "  +11.56%  "
Set-TextValue "D43" "1.998.68"
Set-TextValue "E43" "  +8.54%  "
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  +0.07%  "
Set-TextValue "D45" "91.91"
Set-TextValue "E45" "  +0.56%  "
Set-TextValue "D46" "17.29"
Set-TextValue "E46" "  +37.69%  "
Set-TextValue "D47" "1.83"
Set-TextValue "E47" "  +1.83%  "
Set-TextValue "D48" "107.03"
Set-TextValue "E48" "  +12.57%  "
Set-TextValue "D49" "0.199"
Set-TextValue "E49" "  +9.53%  "
Set-TextValue "D50" "8.89"
Set-TextValue "E50" "  +12.29%  "
Set-TextValue "D51" "2.768.46"
Set-TextValue "E51" "  +8.35%  "
